# Auto-generated Excel COM-interop script
$wb = $excel.ActiveWorkbook

# --- Rename "results" sheet to "trait_h2" and add new "global_rg" sheet ---
$wsDescriptions = $wb.Worksheets.Item(1)
$wsTraitH2 = $wb.Worksheets.Item(2)
$wsTraitH2.Name = "trait_h2"
$wsGlobalRg = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTraitH2)
$wsGlobalRg.Name = "global_rg"

# --- column_descriptions sheet ---
$wsDescriptions.Cells.Clear()
$descData = New-Object "object[,]" 21,3
$descData[0,0] = "sheet"
$descData[0,1] = "column_name"
$descData[0,2] = "description"
$descData[1,0] = "trait_h2"
$descData[1,1] = "phen"
$descData[1,2] = "Phenotype"
$descData[2,0] = "trait_h2"
$descData[2,1] = "total_observed_scale_h2"
$descData[2,2] = "Estimated SNP heritability (h2, observed scale)"
$descData[3,0] = "trait_h2"
$descData[3,1] = "total_observed_scale_h2_se"
$descData[3,2] = "Standard error of h2"
$descData[4,0] = "trait_h2"
$descData[4,1] = "lambda_gc"
$descData[4,2] = "Genomic inflation factor (lambda GC); equivalent of median(chi^2)/0.4549, where denominator indicates expected median of the chi-squared distribution with 1 degree of freedom"
$descData[5,0] = "trait_h2"
$descData[5,1] = "mean_chi_2"
$descData[5,2] = "Mean chi-square statistic"
$descData[6,0] = "trait_h2"
$descData[6,1] = "intercept"
$descData[6,2] = "LD score regression intercept"
$descData[7,0] = "trait_h2"
$descData[7,1] = "intercept_se"
$descData[7,2] = "Standard error for LD score regression intercept"
$descData[8,0] = "trait_h2"
$descData[8,1] = "z"
$descData[8,2] = "Heritability Z-score (equivalent of h2/se)"
$descData[9,0] = "global_rg"
$descData[9,1] = "p1"
$descData[9,2] = "Phenotype 1"
$descData[10,0] = "global_rg"
$descData[10,1] = "p2"
$descData[10,2] = "Phenotype 2"
$descData[11,0] = "global_rg"
$descData[11,1] = "rg"
$descData[11,2] = "The estimated genetic correlation"
$descData[12,0] = "global_rg"
$descData[12,1] = "se"
$descData[12,2] = "The bootstrap standard error of the genetic correlation estimate"
$descData[13,0] = "global_rg"
$descData[13,1] = "z"
$descData[13,2] = "The bootstrap standard error of the genetic correlation estimate"
$descData[14,0] = "global_rg"
$descData[14,1] = "p"
$descData[14,2] = "P-value for genetic correlation"
$descData[15,0] = "global_rg"
$descData[15,1] = "h2_obs"
$descData[15,2] = "Estimated SNP heritability (h2, observed scale) of the second phenotype "
$descData[16,0] = "global_rg"
$descData[16,1] = "h2_obs_se"
$descData[16,2] = "Standard error of h2 for phenotype 2"
$descData[17,0] = "global_rg"
$descData[17,1] = "h2_int"
$descData[17,2] = "Single-trait LD score regression intercept for phenotype 2"
$descData[18,0] = "global_rg"
$descData[18,1] = "h2_int_se"
$descData[18,2] = "Standard error for single-trait LD score regression intercept for phenotype 2"
$descData[19,0] = "global_rg"
$descData[19,1] = "gcov_int"
$descData[19,2] = "Estimated genetic covariance between p1 and p2"
$descData[20,0] = "global_rg"
$descData[20,1] = "gcov_int_se"
$descData[20,2] = "Bootstrap standard error of cross-trait LD score regression intercept"
$wsDescriptions.Range("A1:C21").Value = $descData
$wsDescriptions.Range("A1:C1").Font.Bold = $true

# --- trait_h2 sheet ---
$wsTraitH2.Cells.Clear()
$th2Data = New-Object "object[,]" 9,8
$th2Data[0,0] = "phen"
$th2Data[0,1] = "total_observed_scale_h2"
$th2Data[0,2] = "total_observed_scale_h2_se"
$th2Data[0,3] = "lambda_gc"
$th2Data[0,4] = "mean_chi_2"
$th2Data[0,5] = "intercept"
$th2Data[0,6] = "intercept_se"
$th2Data[0,7] = "z"
$th2Data[1,0] = "AD2019"
$th2Data[1,1] = [double]"0.015"
$th2Data[1,2] = [double]"0.0021"
$th2Data[1,3] = [double]"1.0864"
$th2Data[1,4] = [double]"1.1254"
$th2Data[1,5] = [double]"1.0027"
$th2Data[1,6] = [double]"0.0108"
$th2Data[1,7] = [double]"7.14285714285714"
$th2Data[2,0] = "AD2019.Kunkle"
$th2Data[2,1] = [double]"0.0713"
$th2Data[2,2] = [double]"0.0114"
$th2Data[2,3] = [double]"1.0926"
$th2Data[2,4] = [double]"1.118"
$th2Data[2,5] = [double]"1.0302"
$th2Data[2,6] = [double]"0.0084"
$th2Data[2,7] = [double]"6.25438596491228"
$th2Data[3,0] = "BIP2021"
$th2Data[3,1] = [double]"0.0708"
$th2Data[3,2] = [double]"0.0027"
$th2Data[3,3] = [double]"1.453"
$th2Data[3,4] = [double]"1.5926"
$th2Data[3,5] = [double]"1.0247"
$th2Data[3,6] = [double]"0.0089"
$th2Data[3,7] = [double]"26.2222222222222"
$th2Data[4,0] = "LBD2020"
$th2Data[4,1] = [double]"0.1711"
$th2Data[4,2] = [double]"0.0755"
$th2Data[4,3] = [double]"1.0225"
$th2Data[4,4] = [double]"1.0245"
$th2Data[4,5] = [double]"1.002"
$th2Data[4,6] = [double]"0.0071"
$th2Data[4,7] = [double]"2.26622516556291"
$th2Data[5,0] = "MDD2019"
$th2Data[5,1] = [double]"0.0598"
$th2Data[5,2] = [double]"0.0023"
$th2Data[5,3] = [double]"1.453"
$th2Data[5,4] = [double]"1.5893"
$th2Data[5,5] = [double]"1.0017"
$th2Data[5,6] = [double]"0.0098"
$th2Data[5,7] = 26
$th2Data[6,0] = "PD2019.ex23andMe.exUKBB"
$th2Data[6,1] = [double]"0.3062"
$th2Data[6,2] = [double]"0.0275"
$th2Data[6,3] = [double]"1.0679"
$th2Data[6,4] = [double]"1.0928"
$th2Data[6,5] = [double]"0.977"
$th2Data[6,6] = [double]"0.0057"
$th2Data[6,7] = [double]"11.1345454545455"
$th2Data[7,0] = "PD2019.meta5.ex23andMe"
$th2Data[7,1] = [double]"0.0186"
$th2Data[7,2] = [double]"0.0017"
$th2Data[7,3] = [double]"1.0895"
$th2Data[7,4] = [double]"1.136"
$th2Data[7,5] = [double]"0.9838"
$th2Data[7,6] = [double]"0.0065"
$th2Data[7,7] = [double]"10.9411764705882"
$th2Data[8,0] = "SCZ2018"
$th2Data[8,1] = [double]"0.4101"
$th2Data[8,2] = [double]"0.0138"
$th2Data[8,3] = [double]"1.6831"
$th2Data[8,4] = [double]"1.932"
$th2Data[8,5] = [double]"1.0699"
$th2Data[8,6] = [double]"0.0113"
$th2Data[8,7] = [double]"29.7173913043478"
$wsTraitH2.Range("A1:H9").Value = $th2Data
$wsTraitH2.Range("A1:H1").Font.Bold = $true
$wsTraitH2.Activate()
$wsTraitH2.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# --- global_rg sheet ---
$wsGlobalRg.Cells.Clear()
$grgData = New-Object "object[,]" 37,12
$grgData[0,0] = "p1"
$grgData[0,1] = "p2"
$grgData[0,2] = "rg"
$grgData[0,3] = "se"
$grgData[0,4] = "z"
$grgData[0,5] = "p"
$grgData[0,6] = "h2_obs"
$grgData[0,7] = "h2_obs_se"
$grgData[0,8] = "h2_int"
$grgData[0,9] = "h2_int_se"
$grgData[0,10] = "gcov_int"
$grgData[0,11] = "gcov_int_se"
$grgData[1,0] = "AD2019"
$grgData[1,1] = "AD2019"
$grgData[1,2] = [double]"1.0"
$grgData[1,3] = [double]"6.8394e-06"
$grgData[1,4] = [double]"146210.3044"
$grgData[1,5] = [double]"0.0"
$grgData[1,6] = [double]"0.0086"
$grgData[1,7] = [double]"0.005"
$grgData[1,8] = [double]"1.0517"
$grgData[1,9] = [double]"0.0575"
$grgData[1,10] = [double]"1.0517"
$grgData[1,11] = [double]"0.0575"
$grgData[2,0] = "AD2019"
$grgData[2,1] = "BIP2021"
$grgData[2,2] = [double]"0.1266"
$grgData[2,3] = [double]"0.0585"
$grgData[2,4] = [double]"2.1647"
$grgData[2,5] = [double]"0.0304"
$grgData[2,6] = [double]"0.0715"
$grgData[2,7] = [double]"0.0026"
$grgData[2,8] = [double]"1.0204"
$grgData[2,9] = [double]"0.0092"
$grgData[2,10] = [double]"0.0065"
$grgData[2,11] = [double]"0.0054"
$grgData[3,0] = "AD2019"
$grgData[3,1] = "LBD2020"
$grgData[3,2] = [double]"0.3769"
$grgData[3,3] = [double]"0.2345"
$grgData[3,4] = [double]"1.6073"
$grgData[3,5] = [double]"0.108"
$grgData[3,6] = [double]"0.176"
$grgData[3,7] = [double]"0.0798"
$grgData[3,8] = [double]"1.0007"
$grgData[3,9] = [double]"0.008"
$grgData[3,10] = [double]"0.0266"
$grgData[3,11] = [double]"0.0149"
$grgData[4,0] = "AD2019"
$grgData[4,1] = "MDD2019"
$grgData[4,2] = [double]"0.1696"
$grgData[4,3] = [double]"0.0676"
$grgData[4,4] = [double]"2.5077"
$grgData[4,5] = [double]"0.0122"
$grgData[4,6] = [double]"0.061"
$grgData[4,7] = [double]"0.0024"
$grgData[4,8] = [double]"0.9942"
$grgData[4,9] = [double]"0.0105"
$grgData[4,10] = [double]"0.0199"
$grgData[4,11] = [double]"0.0055"
$grgData[5,0] = "AD2019"
$grgData[5,1] = "PD2019.meta5.ex23andMe"
$grgData[5,2] = [double]"0.1973"
$grgData[5,3] = [double]"0.0814"
$grgData[5,4] = [double]"2.4233"
$grgData[5,5] = [double]"0.0154"
$grgData[5,6] = [double]"0.0194"
$grgData[5,7] = [double]"0.002"
$grgData[5,8] = [double]"0.9778"
$grgData[5,9] = [double]"0.0075"
$grgData[5,10] = [double]"0.0121"
$grgData[5,11] = [double]"0.0043"
$grgData[6,0] = "AD2019"
$grgData[6,1] = "SCZ2018"
$grgData[6,2] = [double]"0.1087"
$grgData[6,3] = [double]"0.0486"
$grgData[6,4] = [double]"2.2345"
$grgData[6,5] = [double]"0.0254"
$grgData[6,6] = [double]"0.4232"
$grgData[6,7] = [double]"0.0144"
$grgData[6,8] = [double]"1.0515"
$grgData[6,9] = [double]"0.012"
$grgData[6,10] = [double]"0.0121"
$grgData[6,11] = [double]"0.0061"
$grgData[7,0] = "BIP2021"
$grgData[7,1] = "AD2019"
$grgData[7,2] = [double]"0.1266"
$grgData[7,3] = [double]"0.0585"
$grgData[7,4] = [double]"2.1647"
$grgData[7,5] = [double]"0.0304"
$grgData[7,6] = [double]"0.0082"
$grgData[7,7] = [double]"0.0054"
$grgData[7,8] = [double]"1.0556"
$grgData[7,9] = [double]"0.0611"
$grgData[7,10] = [double]"0.0065"
$grgData[7,11] = [double]"0.0054"
$grgData[8,0] = "BIP2021"
$grgData[8,1] = "BIP2021"
$grgData[8,2] = [double]"1.0"
$grgData[8,3] = [double]"6.2656e-07"
$grgData[8,4] = [double]"1596030.0952"
$grgData[8,5] = [double]"0.0"
$grgData[8,6] = [double]"0.0716"
$grgData[8,7] = [double]"0.0028"
$grgData[8,8] = [double]"1.02"
$grgData[8,9] = [double]"0.0091"
$grgData[8,10] = [double]"1.02"
$grgData[8,11] = [double]"0.0091"
$grgData[9,0] = "BIP2021"
$grgData[9,1] = "LBD2020"
$grgData[9,2] = [double]"-0.1158"
$grgData[9,3] = [double]"0.0989"
$grgData[9,4] = [double]"-1.1713"
$grgData[9,5] = [double]"0.2415"
$grgData[9,6] = [double]"0.133"
$grgData[9,7] = [double]"0.0797"
$grgData[9,8] = [double]"1.0066"
$grgData[9,9] = [double]"0.0082"
$grgData[9,10] = [double]"0.0081"
$grgData[9,11] = [double]"0.0054"
$grgData[10,0] = "BIP2021"
$grgData[10,1] = "MDD2019"
$grgData[10,2] = [double]"0.4556"
$grgData[10,3] = [double]"0.0217"
$grgData[10,4] = [double]"21.0209"
$grgData[10,5] = [double]"4.2217e-98"
$grgData[10,6] = [double]"0.0607"
$grgData[10,7] = [double]"0.0024"
$grgData[10,8] = [double]"0.9973"
$grgData[10,9] = [double]"0.0101"
$grgData[10,10] = [double]"0.0604"
$grgData[10,11] = [double]"0.0064"
$grgData[11,0] = "BIP2021"
$grgData[11,1] = "PD2019.meta5.ex23andMe"
$grgData[11,2] = [double]"0.0576"
$grgData[11,3] = [double]"0.0342"
$grgData[11,4] = [double]"1.6866"
$grgData[11,5] = [double]"0.0917"
$grgData[11,6] = [double]"0.0192"
$grgData[11,7] = [double]"0.002"
$grgData[11,8] = [double]"0.9787"
$grgData[11,9] = [double]"0.0074"
$grgData[11,10] = [double]"-0.0015"
$grgData[11,11] = [double]"0.0057"
$grgData[12,0] = "BIP2021"
$grgData[12,1] = "SCZ2018"
$grgData[12,2] = [double]"0.6925"
$grgData[12,3] = [double]"0.0174"
$grgData[12,4] = [double]"39.7256"
$grgData[12,5] = [double]"0.0"
$grgData[12,6] = [double]"0.4209"
$grgData[12,7] = [double]"0.0155"
$grgData[12,8] = [double]"1.056"
$grgData[12,9] = [double]"0.0124"
$grgData[12,10] = [double]"0.134"
$grgData[12,11] = [double]"0.0078"
$grgData[13,0] = "LBD2020"
$grgData[13,1] = "AD2019"
$grgData[13,2] = [double]"0.3769"
$grgData[13,3] = [double]"0.2345"
$grgData[13,4] = [double]"1.6073"
$grgData[13,5] = [double]"0.108"
$grgData[13,6] = [double]"0.008"
$grgData[13,7] = [double]"0.0054"
$grgData[13,8] = [double]"1.0564"
$grgData[13,9] = [double]"0.0621"
$grgData[13,10] = [double]"0.0266"
$grgData[13,11] = [double]"0.0149"
$grgData[14,0] = "LBD2020"
$grgData[14,1] = "BIP2021"
$grgData[14,2] = [double]"-0.1158"
$grgData[14,3] = [double]"0.0989"
$grgData[14,4] = [double]"-1.1713"
$grgData[14,5] = [double]"0.2415"
$grgData[14,6] = [double]"0.0723"
$grgData[14,7] = [double]"0.0027"
$grgData[14,8] = [double]"1.0205"
$grgData[14,9] = [double]"0.0093"
$grgData[14,10] = [double]"0.0081"
$grgData[14,11] = [double]"0.0054"
$grgData[15,0] = "LBD2020"
$grgData[15,1] = "LBD2020"
$grgData[15,2] = [double]"1.0"
$grgData[15,3] = [double]"8.9832e-06"
$grgData[15,4] = [double]"111318.1449"
$grgData[15,5] = [double]"0.0"
$grgData[15,6] = [double]"0.156"
$grgData[15,7] = [double]"0.0807"
$grgData[15,8] = [double]"1.0039"
$grgData[15,9] = [double]"0.0082"
$grgData[15,10] = [double]"1.0039"
$grgData[15,11] = [double]"0.0082"
$grgData[16,0] = "LBD2020"
$grgData[16,1] = "MDD2019"
$grgData[16,2] = [double]"-0.004"
$grgData[16,3] = [double]"0.084"
$grgData[16,4] = [double]"-0.0478"
$grgData[16,5] = [double]"0.9619"
$grgData[16,6] = [double]"0.0614"
$grgData[16,7] = [double]"0.0025"
$grgData[16,8] = [double]"0.9917"
$grgData[16,9] = [double]"0.0106"
$grgData[16,10] = [double]"0.0004"
$grgData[16,11] = [double]"0.0054"
$grgData[17,0] = "LBD2020"
$grgData[17,1] = "PD2019.meta5.ex23andMe"
$grgData[17,2] = [double]"0.6238"
$grgData[17,3] = [double]"0.1692"
$grgData[17,4] = [double]"3.6869"
$grgData[17,5] = [double]"0.0002"
$grgData[17,6] = [double]"0.0193"
$grgData[17,7] = [double]"0.002"
$grgData[17,8] = [double]"0.9797"
$grgData[17,9] = [double]"0.0076"
$grgData[17,10] = [double]"0.0127"
$grgData[17,11] = [double]"0.0046"
$grgData[18,0] = "LBD2020"
$grgData[18,1] = "SCZ2018"
$grgData[18,2] = [double]"-0.0312"
$grgData[18,3] = [double]"0.0755"
$grgData[18,4] = [double]"-0.4134"
$grgData[18,5] = [double]"0.6793"
$grgData[18,6] = [double]"0.4211"
$grgData[18,7] = [double]"0.0144"
$grgData[18,8] = [double]"1.0543"
$grgData[18,9] = [double]"0.0114"
$grgData[18,10] = [double]"0.0068"
$grgData[18,11] = [double]"0.0061"
$grgData[19,0] = "MDD2019"
$grgData[19,1] = "AD2019"
$grgData[19,2] = [double]"0.1696"
$grgData[19,3] = [double]"0.0676"
$grgData[19,4] = [double]"2.5077"
$grgData[19,5] = [double]"0.0122"
$grgData[19,6] = [double]"0.0084"
$grgData[19,7] = [double]"0.0052"
$grgData[19,8] = [double]"1.0525"
$grgData[19,9] = [double]"0.0592"
$grgData[19,10] = [double]"0.0199"
$grgData[19,11] = [double]"0.0055"
$grgData[20,0] = "MDD2019"
$grgData[20,1] = "BIP2021"
$grgData[20,2] = [double]"0.4556"
$grgData[20,3] = [double]"0.0217"
$grgData[20,4] = [double]"21.0209"
$grgData[20,5] = [double]"4.2217e-98"
$grgData[20,6] = [double]"0.0719"
$grgData[20,7] = [double]"0.0026"
$grgData[20,8] = [double]"1.0176"
$grgData[20,9] = [double]"0.009"
$grgData[20,10] = [double]"0.0604"
$grgData[20,11] = [double]"0.0064"
$grgData[21,0] = "MDD2019"
$grgData[21,1] = "LBD2020"
$grgData[21,2] = [double]"-0.004"
$grgData[21,3] = [double]"0.084"
$grgData[21,4] = [double]"-0.0478"
$grgData[21,5] = [double]"0.9619"
$grgData[21,6] = [double]"0.1828"
$grgData[21,7] = [double]"0.0825"
$grgData[21,8] = [double]"1.0"
$grgData[21,9] = [double]"0.0087"
$grgData[21,10] = [double]"0.0004"
$grgData[21,11] = [double]"0.0054"
$grgData[22,0] = "MDD2019"
$grgData[22,1] = "MDD2019"
$grgData[22,2] = [double]"1.0"
$grgData[22,3] = [double]"3.4462e-09"
$grgData[22,4] = [double]"290170000.0"
$grgData[22,5] = [double]"0.0"
$grgData[22,6] = [double]"0.0608"
$grgData[22,7] = [double]"0.0024"
$grgData[22,8] = [double]"0.9946"
$grgData[22,9] = [double]"0.0103"
$grgData[22,10] = [double]"0.9946"
$grgData[22,11] = [double]"0.0103"
$grgData[23,0] = "MDD2019"
$grgData[23,1] = "PD2019.meta5.ex23andMe"
$grgData[23,2] = [double]"-0.0135"
$grgData[23,3] = [double]"0.0334"
$grgData[23,4] = [double]"-0.4029"
$grgData[23,5] = [double]"0.687"
$grgData[23,6] = [double]"0.0192"
$grgData[23,7] = [double]"0.002"
$grgData[23,8] = [double]"0.979"
$grgData[23,9] = [double]"0.0075"
$grgData[23,10] = [double]"0.0043"
$grgData[23,11] = [double]"0.0056"
$grgData[24,0] = "MDD2019"
$grgData[24,1] = "SCZ2018"
$grgData[24,2] = [double]"0.3289"
$grgData[24,3] = [double]"0.0216"
$grgData[24,4] = [double]"15.236"
$grgData[24,5] = [double]"2.0393e-52"
$grgData[24,6] = [double]"0.4214"
$grgData[24,7] = [double]"0.0147"
$grgData[24,8] = [double]"1.0546"
$grgData[24,9] = [double]"0.0118"
$grgData[24,10] = [double]"0.037"
$grgData[24,11] = [double]"0.0078"
$grgData[25,0] = "PD2019.meta5.ex23andMe"
$grgData[25,1] = "AD2019"
$grgData[25,2] = [double]"0.1973"
$grgData[25,3] = [double]"0.0814"
$grgData[25,4] = [double]"2.4233"
$grgData[25,5] = [double]"0.0154"
$grgData[25,6] = [double]"0.0083"
$grgData[25,7] = [double]"0.0051"
$grgData[25,8] = [double]"1.054"
$grgData[25,9] = [double]"0.0588"
$grgData[25,10] = [double]"0.0121"
$grgData[25,11] = [double]"0.0043"
$grgData[26,0] = "PD2019.meta5.ex23andMe"
$grgData[26,1] = "BIP2021"
$grgData[26,2] = [double]"0.0576"
$grgData[26,3] = [double]"0.0342"
$grgData[26,4] = [double]"1.6866"
$grgData[26,5] = [double]"0.0917"
$grgData[26,6] = [double]"0.0717"
$grgData[26,7] = [double]"0.0026"
$grgData[26,8] = [double]"1.0194"
$grgData[26,9] = [double]"0.0091"
$grgData[26,10] = [double]"-0.0015"
$grgData[26,11] = [double]"0.0057"
$grgData[27,0] = "PD2019.meta5.ex23andMe"
$grgData[27,1] = "LBD2020"
$grgData[27,2] = [double]"0.6238"
$grgData[27,3] = [double]"0.1692"
$grgData[27,4] = [double]"3.6869"
$grgData[27,5] = [double]"0.0002"
$grgData[27,6] = [double]"0.1697"
$grgData[27,7] = [double]"0.0833"
$grgData[27,8] = [double]"1.0018"
$grgData[27,9] = [double]"0.0087"
$grgData[27,10] = [double]"0.0127"
$grgData[27,11] = [double]"0.0046"
$grgData[28,0] = "PD2019.meta5.ex23andMe"
$grgData[28,1] = "MDD2019"
$grgData[28,2] = [double]"-0.0135"
$grgData[28,3] = [double]"0.0334"
$grgData[28,4] = [double]"-0.4029"
$grgData[28,5] = [double]"0.687"
$grgData[28,6] = [double]"0.0607"
$grgData[28,7] = [double]"0.0024"
$grgData[28,8] = [double]"0.9956"
$grgData[28,9] = [double]"0.0102"
$grgData[28,10] = [double]"0.0043"
$grgData[28,11] = [double]"0.0056"
$grgData[29,0] = "PD2019.meta5.ex23andMe"
$grgData[29,1] = "PD2019.meta5.ex23andMe"
$grgData[29,2] = [double]"1.0"
$grgData[29,3] = [double]"1.1616e-07"
$grgData[29,4] = [double]"8608520.3621"
$grgData[29,5] = [double]"0.0"
$grgData[29,6] = [double]"0.0193"
$grgData[29,7] = [double]"0.0019"
$grgData[29,8] = [double]"0.9787"
$grgData[29,9] = [double]"0.0073"
$grgData[29,10] = [double]"0.9787"
$grgData[29,11] = [double]"0.0073"
$grgData[30,0] = "PD2019.meta5.ex23andMe"
$grgData[30,1] = "SCZ2018"
$grgData[30,2] = [double]"0.0239"
$grgData[30,3] = [double]"0.0307"
$grgData[30,4] = [double]"0.7777"
$grgData[30,5] = [double]"0.4367"
$grgData[30,6] = [double]"0.4214"
$grgData[30,7] = [double]"0.0145"
$grgData[30,8] = [double]"1.0544"
$grgData[30,9] = [double]"0.0116"
$grgData[30,10] = [double]"0.003"
$grgData[30,11] = [double]"0.0058"
$grgData[31,0] = "SCZ2018"
$grgData[31,1] = "AD2019"
$grgData[31,2] = [double]"0.1087"
$grgData[31,3] = [double]"0.0486"
$grgData[31,4] = [double]"2.2345"
$grgData[31,5] = [double]"0.0254"
$grgData[31,6] = [double]"0.0083"
$grgData[31,7] = [double]"0.0054"
$grgData[31,8] = [double]"1.0562"
$grgData[31,9] = [double]"0.0624"
$grgData[31,10] = [double]"0.0121"
$grgData[31,11] = [double]"0.0061"
$grgData[32,0] = "SCZ2018"
$grgData[32,1] = "BIP2021"
$grgData[32,2] = [double]"0.6925"
$grgData[32,3] = [double]"0.0174"
$grgData[32,4] = [double]"39.7256"
$grgData[32,5] = [double]"0.0"
$grgData[32,6] = [double]"0.0718"
$grgData[32,7] = [double]"0.0028"
$grgData[32,8] = [double]"1.0178"
$grgData[32,9] = [double]"0.0096"
$grgData[32,10] = [double]"0.134"
$grgData[32,11] = [double]"0.0078"
$grgData[33,0] = "SCZ2018"
$grgData[33,1] = "LBD2020"
$grgData[33,2] = [double]"-0.0312"
$grgData[33,3] = [double]"0.0755"
$grgData[33,4] = [double]"-0.4134"
$grgData[33,5] = [double]"0.6793"
$grgData[33,6] = [double]"0.16"
$grgData[33,7] = [double]"0.0827"
$grgData[33,8] = [double]"1.0031"
$grgData[33,9] = [double]"0.0084"
$grgData[33,10] = [double]"0.0068"
$grgData[33,11] = [double]"0.0061"
$grgData[34,0] = "SCZ2018"
$grgData[34,1] = "MDD2019"
$grgData[34,2] = [double]"0.3289"
$grgData[34,3] = [double]"0.0216"
$grgData[34,4] = [double]"15.236"
$grgData[34,5] = [double]"2.0393e-52"
$grgData[34,6] = [double]"0.0605"
$grgData[34,7] = [double]"0.0025"
$grgData[34,8] = [double]"0.9985"
$grgData[34,9] = [double]"0.0103"
$grgData[34,10] = [double]"0.037"
$grgData[34,11] = [double]"0.0078"
$grgData[35,0] = "SCZ2018"
$grgData[35,1] = "PD2019.meta5.ex23andMe"
$grgData[35,2] = [double]"0.0239"
$grgData[35,3] = [double]"0.0307"
$grgData[35,4] = [double]"0.7777"
$grgData[35,5] = [double]"0.4367"
$grgData[35,6] = [double]"0.019"
$grgData[35,7] = [double]"0.0018"
$grgData[35,8] = [double]"0.9796"
$grgData[35,9] = [double]"0.0073"
$grgData[35,10] = [double]"0.003"
$grgData[35,11] = [double]"0.0058"
$grgData[36,0] = "SCZ2018"
$grgData[36,1] = "SCZ2018"
$grgData[36,2] = [double]"1.0"
$grgData[36,3] = [double]"1.8166e-08"
$grgData[36,4] = [double]"55047950.1835"
$grgData[36,5] = [double]"0.0"
$grgData[36,6] = [double]"0.4217"
$grgData[36,7] = [double]"0.015"
$grgData[36,8] = [double]"1.0533"
$grgData[36,9] = [double]"0.0121"
$grgData[36,10] = [double]"1.0533"
$grgData[36,11] = [double]"0.0121"
$wsGlobalRg.Range("A1:L37").Value = $grgData
$wsGlobalRg.Range("A1:L1").Font.Bold = $true
$wsGlobalRg.Activate()
$wsGlobalRg.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$wsDescriptions.Activate()
Write-Host "done"
